# additional F_min expressions tested
$wb = $excel.ActiveWorkbook

# --- constants_evaluated sheet: refreshed F_min constant values ---
$wsConst = $wb.Worksheets.Item("constants_evaluated")
$wsConst.Range("C3").Value = 0.00005575661819269246
$wsConst.Range("C4").Value = 0.00006103324607390475

# --- enthalpies_calc sheet: row 2 & row 3 updated values ---
$wsEnth = $wb.Worksheets.Item("enthalpies_calc")
$wsEnth.Range("B2").Value = -0.0
$wsEnth.Range("C2").Value = -0.0
$wsEnth.Range("D2").Value = -9.99999886117819
$wsEnth.Range("E2").Value = -13.00008072751907
$wsEnth.Range("D3").Value = 0.000002842613907042811
$wsEnth.Range("E3").Value = 0.00002549392346979301

# --- correlation_matrix sheet: row 2 & row 3 updated values ---
$wsCorr = $wb.Worksheets.Item("correlation_matrix")
$wsCorr.Range("B2").Value = 1
$wsCorr.Range("C2").Value = -0.9995575693907462
$wsCorr.Range("B3").Value = -0.9995575693907465
